$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "4209-RBI-EI-DB-DL-REC-RNI-FEE-FFConMONTHLYonDAY25-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-ONT-PER-1st"

# Update product name on both sheets (B1)
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# Update short name on input sheet (B2) - now text instead of a number
$ws1.Range("B2").Value = "420q"

# Reset selection on the input sheet back to the top-left cell
$ws1.Range("B1").Select()

# Make the output sheet the active / selected tab
$ws2.Activate()
